$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before old row 2 (pushes old rows 2-23 down to 6-27)
$ws.Rows("2:5").Insert()

# New row 2: section header for the new SQL-database test case
$ws.Range("A2").Value = "New tables in SQL database are created"

# New row 3: new test-case data row
$ws.Range("A3").Value = "New track model has 3 lines"
$ws.Range("B3").Value = "1) Select track model file`n2) Click import"
$ws.Range("C3").Value = "3 tables (Blue Line, Red Line, and Green Line) are created in the SQL database"
$ws.Range("D3").Value = "Pass"
$ws.Range("F3").Value = 44893
$ws.Range("F3").NumberFormat = "d-mmm"

$ws.Range("A3:D3").Font.Bold = $false
$ws.Range("A3:D3").WrapText = $true
$ws.Rows("3:3").RowHeight = 44.25

# Row 7 (old row 3, shifted down) keeps its data but loses its bold/header styling
$ws.Range("A7:C7").WrapText = $true

# Column A is widened to fit the new, longer test-case text
$ws.Columns("A:A").ColumnWidth = 39.5

# Update the saved view: no frozen/scrolled top row, selection on B4
$null = $ws.Range("B4").Select()
